# evaluateRegModel.xlsx: refresh the "x" feature-list dictionary (row 2,
# column B on the Indicators sheet) with the new set of encoded predictor
# columns, and widen column B so the longer string still fits/best-fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicators")

$ws.Range("B2").Value = "SchoolRegion_2, SchoolType_2, MotherEd_4, MotherEd_7, PostulationType_1"

# Column B held a bestFit width for the old, shorter string (51.4 chars).
# Re-fit it to the new, longer string.
$ws.Columns("B").ColumnWidth = 73.333333333
